$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data as scraped on Wed Aug  7 16:45:32 UTC 2024
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.018.76'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.387.07'
$ws.Range('E3').Value = '  -4.75%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '479.05'
$ws.Range('E5').Value = '  -2.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.91'
$ws.Range('E6').Value = '  +1.15%  '
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('E8').Value = '  -2.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.383.23'
$ws.Range('E9').Value = '  -5.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0974'
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.44'
$ws.Range('E11').Value = '  -5.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.323'
$ws.Range('E12').Value = '  -3.20%  '
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.802.98'
$ws.Range('E14').Value = '  -4.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '56.091.72'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.32'
$ws.Range('E16').Value = '  -4.26%  '
$ws.Range('E17').Value = '  -3.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.397.52'
$ws.Range('E18').Value = '  -4.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.52'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '314.72'
$ws.Range('E20').Value = '  -2.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.71'
$ws.Range('E21').Value = '  -5.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.68'
$ws.Range('E23').Value = '  -2.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '56.83'
$ws.Range('E24').Value = '  -3.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.394'
$ws.Range('E26').Value = '  -4.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.158'
$ws.Range('E27').Value = '  -5.57%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.499.65'
$ws.Range('E28').Value = '  -4.45%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.24'
$ws.Range('E29').Value = '  -4.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0770'
$ws.Range('E30').Value = '  -3.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.02'
$ws.Range('E32').Value = '  -0.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.97'
$ws.Range('E33').Value = '  -2.25%  '
$ws.Range('E34').Value = '  -1.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.03'
$ws.Range('E35').Value = '  -3.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.10'
$ws.Range('E36').Value = '  -4.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.58'
$ws.Range('E37').Value = '  -4.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.837'
$ws.Range('E38').Value = '  -3.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '33.40'
$ws.Range('E39').Value = '  -2.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  +0.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.34'
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.38'
$ws.Range('E42').Value = '  -4.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0537'
$ws.Range('E43').Value = '  -3.67%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0943'
$ws.Range('E44').Value = '  +2.96%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.582'
$ws.Range('E45').Value = '  -5.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.23'
$ws.Range('E46').Value = '  +0.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '254.46'
$ws.Range('E47').Value = '  -2.62%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.61'
$ws.Range('E48').Value = '  -5.36%  '
$ws.Range('E49').Value = '  -2.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.99'
$ws.Range('E50').Value = '  -3.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.783.50'
$ws.Range('E51').Value = '  -7.37%  '
